$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.053.09'
$ws.Range("E2").Value = '  +0.43%  '
$ws.Range("D3").Value = '3.207.11'
$ws.Range("E3").Value = '  -0.78%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.68'
$ws.Range("E5").Value = '  -0.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.22'
$ws.Range("E6").Value = '  -3.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.597'
$ws.Range("E7").Value = '  -5.16%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  -2.60%  '
$ws.Range("E10").Value = '  -0.19%  '
$ws.Range("D12").Value = '3.762.30'
$ws.Range("E12").Value = '  -0.72%  '
$ws.Range("D14").Value = '65.065.86'
$ws.Range("E14").Value = '  +0.29%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.57'
$ws.Range("E15").Value = '  -0.77%  '
$ws.Range("D16").Value = '3.209.73'
$ws.Range("E16").Value = '  -0.67%  '
$ws.Range("E17").Value = '  -0.96%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '411.66'
$ws.Range("E18").Value = '  -1.24%  '
$ws.Range("E19").Value = '  +0.85%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.34'
$ws.Range("E20").Value = '  -0.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.19'
$ws.Range("E21").Value = '  -0.34%  '
$ws.Range("E22").Value = '  +0.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.68'
$ws.Range("E23").Value = '  -1.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.202'
$ws.Range("E24").Value = '  -1.58%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.492'
$ws.Range("E25").Value = '  -0.67%  '
$ws.Range("E26").Value = '  -5.72%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.90'
$ws.Range("E27").Value = '  -1.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.997'
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("E29").Value = '  -1.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '21.55'
$ws.Range("E30").Value = '  -1.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.97'
$ws.Range("E31").Value = '  -0.23%  '
$ws.Range("E32").Value = '  -0.21%  '
$ws.Range("E33").Value = '  -1.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '156.86'
$ws.Range("E34").Value = '  -0.62%  '
$ws.Range("E35").Value = '  -1.52%  '
$ws.Range("B36").Value = 'Maker'
$ws.Range("C36").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D36").Value = '2.755.09'
$ws.Range("E36").Value = '  -2.51%  '
$ws.Range("B37").Value = 'Stacks'
$ws.Range("C37").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.73'
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '24.26'
$ws.Range("E38").Value = '  -4.52%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.15'
$ws.Range("E39").Value = '  -1.67%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.714'
$ws.Range("E40").Value = '  -1.37%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0633'
$ws.Range("E41").Value = '  +0.62%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.65'
$ws.Range("E42").Value = '  -1.68%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0263'
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '296.30'
$ws.Range("E44").Value = '  -1.75%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.49'
$ws.Range("E45").Value = '  -2.29%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").Value = '  -0.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0992'
$ws.Range("E47").Value = '  -1.86%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.96'
$ws.Range("E48").Value = '  -9.60%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.80'
$ws.Range("E49").Value = '  -0.39%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '10.46'
$ws.Range("E50").Value = '  +0.58%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.903'
$ws.Range("E51").Value = '  -2.71%  '
